{"js": "const table = context.document.body.tables.getFirst();\n\ntable.getCell(0, 0).value = \"0M\";\ntable.getCell(1, 0).value = \"0M\";\ntable.getCell(2, 0).value = \"0M\";\ntable.getCell(3, 0).value = \"75\";\n\ntable.getCell(5, 0).value = \"0.00074\";\ntable.getCell(6, 0).value = \"0.00018\";\ntable.getCell(7, 0).value = \"0.00005\";\ntable.getCell(8, 0).value = \"0.00030\";\ntable.getCell(9, 0).value = \"0.00040\";\ntable.getCell(10, 0).value = \"0.00044\";\ntable.getCell(11, 0).value = \"0.01615\";\n\ntable.getCell(43, 0).value = \"99.97\";\ntable.getCell(44, 0).value = \"0.02\";\ntable.getCell(45, 0).value = \"56\";\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$t.Cell(1, 1).Range.Text = \"0M\"\n$t.Cell(2, 1).Range.Text = \"0M\"\n$t.Cell(3, 1).Range.Text = \"0M\"\n$t.Cell(4, 1).Range.Text = \"75\"\n\n$t.Cell(6, 1).Range.Text = \"0.00074\"\n$t.Cell(7, 1).Range.Text = \"0.00018\"\n$t.Cell(8, 1).Range.Text = \"0.00005\"\n$t.Cell(9, 1).Range.Text = \"0.00030\"\n$t.Cell(10, 1).Range.Text = \"0.00040\"\n$t.Cell(11, 1).Range.Text = \"0.00044\"\n$t.Cell(12, 1).Range.Text = \"0.01615\"\n\n$t.Cell(44, 1).Range.Text = \"99.97\"\n$t.Cell(45, 1).Range.Text = \"0.02\"\n$t.Cell(46, 1).Range.Text = \"56\"\n"}
